# Edit script for "10 QUALITY MANAGEMENT PLAN - DEVELOPMENTALITY.docx"
# Applies the textual changes described by the diff:
#   1. Merge the "This section describes..." paragraph's three runs (which
#      had proofing marks around "in order to") into plain continuous text.
#   2. Replace the ": " separator after several bullet headings
#      (Code Review, Unit Testing, Integration Testing,
#       User Acceptance Testing (UAT), Security Testing,
#       Continuous Monitoring and Maintenance) with " - ".
#   3. Remove stray proofing-mark line breaks in the quality-control bullet
#      list by collapsing their two runs into continuous text
#      (assessing .. results / corrective .. measures / remedial .. completed
#       / implementing corrective .. measures).

$d = $word.ActiveDocument

# 1) Merge the intro paragraph text (drop the "in order to" proofing split).
$d.Content.Find.Execute(
    "This section describes the approach the organization will use for managing quality throughout the project's life cycle.  Quality must always be planned into a project in order to prevent unnecessary rework, waste, cost, and time.  Quality should also be considered from both a product and process perspective.  The organization may already have a standardized approach to quality, however, whether it is standard or not, the approach must be defined and communicated to all project stakeholders.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This section describes the approach the organization will use for managing quality throughout the project’s life cycle.  Quality must always be planned into a project in order to prevent unnecessary rework, waste, cost, and time.  Quality should also be considered from both a product and process perspective.  The organization may already have a standardized approach to quality, however, whether it is standard or not, the approach must be defined and communicated to all project stakeholders.",
    2)

# 2) Replace ": " with " - " after each of these headings.
$headings = @(
    "Code Review",
    "Unit Testing",
    "Integration Testing",
    "User Acceptance Testing (UAT)",
    "Security Testing",
    "Continuous Monitoring and Maintenance"
)

foreach ($h in $headings) {
    $rng = $d.Content
    $rng.Find.Execute($h, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($rng.Find.Found) {
        # Move to just after the heading text, then grab the next couple
        # characters (the ": ") and replace them with " - ".
        $afterRng = $d.Range($rng.End, $rng.End + 2)
        if ($afterRng.Text -eq ": ") {
            $afterRng.Text = " - "
        }
    }
}

# 3) Collapse the split bullet-list sentences (remove the mid-sentence
#    proofing-mark break) by re-joining the two runs of text.
$d.Content.Find.Execute(
    "Team member responsible for assessing the measurement results",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Team member responsible for assessing the measurement results", 2)

$d.Content.Find.Execute(
    "Actions taken for any required corrective measures",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Actions taken for any required corrective measures", 2)

$d.Content.Find.Execute(
    "Date when the remedial measures were completed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Date when the remedial measures were completed", 2)

$d.Content.Find.Execute(
    "Team member responsible for implementing corrective measures",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Team member responsible for implementing corrective measures", 2)
